$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for Price (D) cells whose new values look numeric,
# so Excel keeps them as text (preserving trailing zeros / leading formatting)
# instead of silently converting to a number, matching the source data.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.213.37"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "3.783.12"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "594.77"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("D6").Value = "171.05"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("D7").Value = "3.776.31"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").Value = "6.25"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "38.26"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "4.427.48"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").Value = "3.810.42"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "68.619.19"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.116"
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").Value = "15.98"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").Value = "488.28"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").Value = "86.30"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -5.64%  "
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("D27").Value = "12.29"
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  -10.36%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "2.43"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "32.19"
$ws.Range("E32").Value = "  +7.21%  "
$ws.Range("D33").Value = "7.64"
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Value = "5.85"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").Value = "0.135"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").Value = "0.324"
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("D40").Value = "450.14"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("D41").Value = "49.19"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "2.01"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "2.82"
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("D44").Value = "8.33"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "41.29"
$ws.Range("E45").Value = "  -5.94%  "
$ws.Range("D46").Value = "2.841.37"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").Value = "0.0353"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "138.80"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "26.76"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -3.75%  "
